# Applies the "Add files via upload" edit to Saldo_guide.xlsx:
#  - rename the sheet (new export timestamp)
#  - bump the reference date (column G) for every data row from 45589 to 45590
#  - update a handful of rows whose projected/forecast balances (D, E, H) changed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new export timestamp.
$ws.Name = "IClientBalance-20241025-093204-"

# Column G ("Dt. Referencia") moves from 45589 to 45590 for every data row (2-274).
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value2 = 45590
}

# Rows whose balances were recalculated also need their D / E / H values updated.
$rowUpdates = @{
    5   = @{ E = 999.99;               H = 999.99 }
    15  = @{ D = 0;                    E = 8362.76;              H = 8362.76 }
    51  = @{ E = 18550.43;             H = 18550.43 }
    57  = @{ E = 42862.8;              H = 42862.8 }
    60  = @{ E = 999.99;               H = 999.99 }
    104 = @{ E = 999.99;               H = 999.99 }
    105 = @{ E = 4935.71;              H = 4935.71 }
    109 = @{ E = 16527.990000000002;   H = 16527.990000000002 }
    143 = @{ E = 69709.69;             H = 69709.69 }
    173 = @{ E = 999.99;               H = 999.99 }
    232 = @{ E = 44941.99;             H = 44941.99 }
    235 = @{ E = 956.72;               H = 956.72 }
    264 = @{ D = -288.2;               E = 1227.43;              H = 939.23 }
    265 = @{ E = 999.99;               H = 999.99 }
    270 = @{ E = 999.99;               H = 999.99 }
    271 = @{ E = 999.99;               H = 999.99 }
    273 = @{ E = 999.99;               H = 999.99 }
}

$colIndex = @{ D = 4; E = 5; H = 8 }

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Cells.Item([int]$row, $colIndex[$col]).Value2 = $cols[$col]
    }
}
